$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued numeric-looking Price cells keep their exact literal
# formatting (trailing zeros / scientific-notation-prone values) by marking
# the target cell as Text before assigning, matching source data which is
# stored as plain strings, not numbers.

$ws.Range("D2").Value = '27.130.82'
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").Value = '1.799.64'
$ws.Range("E3").Value = '  -2.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.006'
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.30'
$ws.Range("E6").Value = '  -1.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4184'
$ws.Range("E7").Value = '  -1.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3563'
$ws.Range("E8").Value = '  -2.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07078'
$ws.Range("E9").Value = '  -2.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8454'
$ws.Range("E10").Value = '  -2.79%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.21'
$ws.Range("E11").Value = '  -2.74%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.783.21'
$ws.Range("E12").Value = '  -2.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.276'
$ws.Range("E13").Value = '  -2.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.351'
$ws.Range("E14").Value = '  -2.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06795'
$ws.Range("E15").Value = '  -2.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.009'
$ws.Range("E16").Value = '  +0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '79.96'
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008719'
$ws.Range("E18").Value = '  -3.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.006'
$ws.Range("E19").Value = '  +0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.10'
$ws.Range("E20").Value = '  -2.40%  '

$ws.Range("D21").Value = '28.098.54'
$ws.Range("E21").Value = '  +2.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.058'
$ws.Range("E22").Value = '  -0.15%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.79'
$ws.Range("E23").Value = '  -0.68%  '

$ws.Range("D24").Value = '2.191.74'
$ws.Range("E24").Value = '  +5.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.959'
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.24'
$ws.Range("E26").Value = '  -0.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.12'
$ws.Range("E27").Value = '  -1.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.045'
$ws.Range("E28").Value = '  -3.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.81'
$ws.Range("E29").Value = '  -2.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.665'
$ws.Range("E30").Value = '  -10.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08878'
$ws.Range("E31").Value = '  +0.02%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7215'
$ws.Range("E32").Value = '  -7.08%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.864'
$ws.Range("E33").Value = '  -3.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.350'
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.006'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.080'
$ws.Range("E36").Value = '  -6.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.075'
$ws.Range("E37").Value = '  -2.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01897'
$ws.Range("E38").Value = '  -2.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05100'
$ws.Range("E39").Value = '  -5.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1619'
$ws.Range("E40").Value = '  -2.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4947'
$ws.Range("E41").Value = '  -3.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.595'
$ws.Range("E42").Value = '  -7.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.214'
$ws.Range("E43").Value = '  -7.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.087'
$ws.Range("E44").Value = '  -5.49%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.25'
$ws.Range("E45").Value = '  -3.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.72'
$ws.Range("E46").Value = '  -1.72%  '

$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  +0.48%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06313'
$ws.Range("E48").Value = '  -3.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4513'
$ws.Range("E49").Value = '  -4.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.589'
$ws.Range("E50").Value = '  -3.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.28'
$ws.Range("E51").Value = '  -3.41%  '
